$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering for the "Filename" column (D2:D81) after the re-shuffle.
$filenames = @(
    "atcoder_ABC132_F.java",
    "codeforces_369_B.java",
    "codeforces_147_A.java",
    "codeforces_110_B.java",
    "atcoder_AGC007_C.java",
    "atcoder_ABC149_B.java",
    "atcoder_ABC172_D.java",
    "codeforces_276_B.java",
    "atcoder_ABC132_A.java",
    "codeforces_79_A.java",
    "codeforces_49_A.java",
    "codeforces_678_B.java",
    "codeforces_373_B.java",
    "codeforces_30_A.java",
    "atcoder_ABC114_C.java",
    "codeforces_569_A.java",
    "codeforces_59_A.java",
    "codeforces_32_B.java",
    "atcoder_ABC178_B.java",
    "codeforces_203_A.java",
    "codeforces_651_A.java",
    "atcoder_ABC124_A.java",
    "codeforces_669_A.java",
    "atcoder_AGC046_B.java",
    "atcoder_AGC025_A.java",
    "codeforces_8_B.java",
    "codeforces_672_A.java",
    "atcoder_ABC070_B.java",
    "codeforces_546_A.java",
    "atcoder_ABC169_D.java",
    "codeforces_306_A.java",
    "codeforces_579_A.java",
    "atcoder_ABC170_A.java",
    "atcoder_ABC174_C.java",
    "atcoder_ABC139_B.java",
    "atcoder_ARC062_B.java",
    "atcoder_ABC153_A.java",
    "atcoder_ABC125_A.java",
    "codeforces_99_A.java",
    "codeforces_459_A.java",
    "codeforces_622_A.java",
    "codeforces_92_A.java",
    "atcoder_ABC043_B.java",
    "atcoder_ABC051_A.java",
    "atcoder_ABC158_B.java",
    "atcoder_ABC178_A.java",
    "atcoder_AGC006_B.java",
    "codeforces_340_A.java",
    "atcoder_ABC143_A.java",
    "atcoder_ABC122_D.java",
    "codeforces_86_A.java",
    "codeforces_678_A.java",
    "codeforces_242_A.java",
    "atcoder_ABC164_A.java",
    "codeforces_514_A.java",
    "atcoder_ABC158_A.java",
    "codeforces_544_B.java",
    "atcoder_ABC127_B.java",
    "codeforces_190_A.java",
    "atcoder_ABC136_B.java",
    "codeforces_189_A.java",
    "atcoder_ABC169_C.java",
    "atcoder_ARC102_C.java",
    "codeforces_581_A.java",
    "atcoder_ABC108_B.java",
    "atcoder_ABC120_C.java",
    "atcoder_AGC046_A.java",
    "atcoder_ABC168_C.java",
    "codeforces_96_B.java",
    "atcoder_ABC151_A.java",
    "atcoder_AGC002_A.java",
    "atcoder_ABC124_C.java",
    "codeforces_171_A.java",
    "codeforces_379_A.java",
    "atcoder_ABC142_A.java",
    "atcoder_ABC149_C.java",
    "codeforces_58_B.java",
    "codeforces_55_A.java",
    "atcoder_ABC042_A.java",
    "codeforces_334_A.java"
)

# New "Impact" value per row (E2:E81). The category boundaries moved: the
# Test Failed block grew (absorbing what used to be most of the Correct
# block), and the single Infinite Loop row moved down near the bottom.
for ($i = 0; $i -lt $filenames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $filenames[$i]

    if ($row -le 6) {
        $impact = "Compilation Error"
    } elseif ($row -le 16) {
        $impact = "Runtime Error"
    } elseif ($row -le 77) {
        $impact = "Test Failed"
    } elseif ($row -eq 78) {
        $impact = "Infinite Loop"
    } else {
        $impact = "Correct"
    }
    $ws.Cells.Item($row, 5).Value = $impact
}
